$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at sheet row 17 (shifts "Affliction's Heart" and everything
# below it down by one row), then grow the structured table to include it.
$ws.Rows.Item(17).Insert()
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E20"))

# Row 16 "Elrunez Confrontation" description gets expanded wording.
$ws.Range("B16").Value = "It's finally come to it. A battle with Elrunez, the Divine Ruler, the Lord of Torment, and the creator of the player. He created the player as a vessel that could create and wield the Atlas of Frosyni, but he seems to have failed. Or so it seems."

# New row 17: "Lord of Torment" - Elrunez's second form.
$ws.Range("A17").Value = "Lord of Torment"
$ws.Range("B17").Value = "This is Elrunez's second form."
$ws.Range("C17").Value = "Elrunez, Affliction"
$ws.Range("E17").Value = "Anguish, Holy"

# Former row 17 ("Affliction's Heart"), now row 18: rename the track.
$ws.Range("A18").Value = "Affliction's Image"

# Row 3 "Insanity" description: Faceless -> Acharos.
$ws.Range("B3").Value = "Plays when the player has low sanity and is attacked by tormentors. Might also play when the player encounters a non-hostile Acharos."

# Row 9 "The Lord's Phantasm": italicize the track-name cell.
$ws.Range("A9").Font.Italic = $true

# Final selection ends up on C3.
[void]$ws.Range("C3").Select()
